$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 290.44446
$ws.Range("I4").Value = 290.44446
$ws.Range("K4").Value = 290.44446
$ws.Range("M4").Value = -176.44446
$ws.Range("H33").Value = 176.6
$ws.Range("I33").Value = 176.6
$ws.Range("K33").Value = 176.6
$ws.Range("M33").Value = 52.40000000000001
$ws.Range("H62").Value = 4082.0476
$ws.Range("J62").Value = 4817.4287
$ws.Range("L62").Value = 4817.4287
$ws.Range("N62").Value = -6065.4287
$ws.Range("H65").Value = 4082.0476
$ws.Range("J65").Value = 4817.4287
$ws.Range("L65").Value = 24087.1435
$ws.Range("N65").Value = -30327.1435
$ws.Range("H70").Value = 2338.4167
$ws.Range("I70").Value = 2408
$ws.Range("J70").Value = 2199.25
$ws.Range("K70").Value = 7224
$ws.Range("L70").Value = 6597.75
$ws.Range("M70").Value = -6954
$ws.Range("N70").Value = -7137.75
$ws.Range("H73").Value = 2338.4167
$ws.Range("I73").Value = 2408
$ws.Range("J73").Value = 2199.25
$ws.Range("K73").Value = 7224
$ws.Range("L73").Value = 6597.75
$ws.Range("M73").Value = -6288
$ws.Range("N73").Value = -8469.75
$ws.Range("H137").Value = 1950.1875
$ws.Range("I137").Value = 1642.0834
$ws.Range("K137").Value = 4926.2502
$ws.Range("M137").Value = -2376.2502
$ws.Range("H138").Value = 2822.5483
$ws.Range("I138").Value = 7576.25
$ws.Range("J138").Value = 2118.2964
$ws.Range("K138").Value = 22728.75
$ws.Range("L138").Value = 6354.889200000001
$ws.Range("M138").Value = -17588.75
$ws.Range("N138").Value = -16634.8892

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 96
$ws.Range("J5").Value = 400
$ws.Range("L5").Value = 400
$ws.Range("N5").Value = -624
$ws.Range("H32").Value = 71197.64999999999
$ws.Range("I32").Value = 13069.77
$ws.Range("K32").Value = 13069.77
$ws.Range("M32").Value = -12782.77
$ws.Range("H45").Value = 6509.5
$ws.Range("I45").Value = 10379.083
$ws.Range("K45").Value = 10379.083
$ws.Range("M45").Value = -10002.083
$ws.Range("H97").Value = 636.6129
$ws.Range("I97").Value = 693.8333
$ws.Range("J97").Value = 440.42856
$ws.Range("K97").Value = 693.8333
$ws.Range("L97").Value = 440.42856
$ws.Range("M97").Value = -197.8333
$ws.Range("N97").Value = -1432.42856
$ws.Range("H110").Value = 1435.8235
$ws.Range("I110").Value = 1359.6666
$ws.Range("K110").Value = 1359.6666
$ws.Range("M110").Value = 685.3334
$ws.Range("H122").Value = 2247.3513
$ws.Range("I122").Value = 2158.8
$ws.Range("K122").Value = 6476.400000000001
$ws.Range("M122").Value = -4026.400000000001
$ws.Range("H132").Value = 2809.2693
$ws.Range("I132").Value = 1552.05
$ws.Range("K132").Value = 4656.15
$ws.Range("M132").Value = -2126.15

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 96
$ws.Range("J4").Value = 400
$ws.Range("L4").Value = 400
$ws.Range("N4").Value = -630
$ws.Range("H76").Value = 20156.75
$ws.Range("I76").Value = 15000
$ws.Range("J76").Value = 21875.666
$ws.Range("K76").Value = 15000
$ws.Range("L76").Value = 21875.666
$ws.Range("M76").Value = -14685
$ws.Range("N76").Value = -22505.666
$ws.Range("H79").Value = 20156.75
$ws.Range("I79").Value = 15000
$ws.Range("J79").Value = 21875.666
$ws.Range("K79").Value = 15000
$ws.Range("L79").Value = 21875.666
$ws.Range("M79").Value = -13908
$ws.Range("N79").Value = -24059.666
$ws.Range("H107").Value = 78306.234
$ws.Range("I107").Value = 101337
$ws.Range("J107").Value = 1537
$ws.Range("K107").Value = 101337
$ws.Range("L107").Value = 1537
$ws.Range("M107").Value = -99417
$ws.Range("N107").Value = -5377

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14208.967
$ws.Range("I31").Value = 3078.1304
$ws.Range("J31").Value = 50781.715
$ws.Range("K31").Value = 3078.1304
$ws.Range("L31").Value = 50781.715
$ws.Range("M31").Value = -2783.1304
$ws.Range("N31").Value = -51371.715
$ws.Range("H34").Value = 14208.967
$ws.Range("I34").Value = 3078.1304
$ws.Range("J34").Value = 50781.715
$ws.Range("K34").Value = 3078.1304
$ws.Range("L34").Value = 50781.715
$ws.Range("M34").Value = -2876.1304
$ws.Range("N34").Value = -51185.715
$ws.Range("H52").Value = 54987
$ws.Range("J52").Value = 59999
$ws.Range("L52").Value = 59999
$ws.Range("N52").Value = -60587
$ws.Range("H74").Value = 80157
$ws.Range("J74").Value = 100314
$ws.Range("L74").Value = 100314
$ws.Range("N74").Value = -102062
$ws.Range("H77").Value = 80157
$ws.Range("J77").Value = 100314
$ws.Range("L77").Value = 300942
$ws.Range("N77").Value = -309678
$ws.Range("H132").Value = 2570.389
$ws.Range("I132").Value = 2508.5806
$ws.Range("K132").Value = 7525.7418
$ws.Range("M132").Value = -4995.7418

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 235
$ws.Range("J7").Value = 70
$ws.Range("L7").Value = 210
$ws.Range("N7").Value = -434
$ws.Range("H14").Value = 1142.875
$ws.Range("I14").Value = 1142.875
$ws.Range("K14").Value = 3428.625
$ws.Range("M14").Value = -3255.625
$ws.Range("H50").Value = 50
$ws.Range("I50").Value = 58.333332
$ws.Range("J50").Value = 25
$ws.Range("K50").Value = 174.999996
$ws.Range("L50").Value = 75
$ws.Range("M50").Value = 306.000004
$ws.Range("N50").Value = -1037
$ws.Range("H53").Value = 50
$ws.Range("I53").Value = 58.333332
$ws.Range("J53").Value = 25
$ws.Range("K53").Value = 174.999996
$ws.Range("L53").Value = 75
$ws.Range("M53").Value = 306.000004
$ws.Range("N53").Value = -1037

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 19511.5
$ws.Range("I58").Value = 17500
$ws.Range("K58").Value = 17500
$ws.Range("M58").Value = -17223
$ws.Range("H93").Value = 40000
$ws.Range("J93").Value = 40000
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -43744
$ws.Range("H97").Value = 34040.57
$ws.Range("I97").Value = 49079.93
$ws.Range("J97").Value = 3961.8572
$ws.Range("K97").Value = 49079.93
$ws.Range("L97").Value = 3961.8572
$ws.Range("M97").Value = -48583.93
$ws.Range("N97").Value = -4953.8572
$ws.Range("H122").Value = 1445.0714
$ws.Range("I122").Value = 1479.4615
$ws.Range("J122").Value = 998
$ws.Range("K122").Value = 4438.3845
$ws.Range("L122").Value = 2994
$ws.Range("M122").Value = -1988.3845
$ws.Range("N122").Value = -7894
$ws.Range("H126").Value = 3237.7585
$ws.Range("I126").Value = 3226.1924
$ws.Range("J126").Value = 3338
$ws.Range("K126").Value = 9678.5772
$ws.Range("L126").Value = 10014
$ws.Range("M126").Value = -7208.5772
$ws.Range("N126").Value = -14954
$ws.Range("H132").Value = 5327.625
$ws.Range("I132").Value = 4626.5415
$ws.Range("J132").Value = 7430.875
$ws.Range("K132").Value = 13879.6245
$ws.Range("L132").Value = 22292.625
$ws.Range("M132").Value = -11349.6245
$ws.Range("N132").Value = -27352.625

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3819.2678
$ws.Range("I132").Value = 3217.2683
$ws.Range("K132").Value = 9651.804900000001
$ws.Range("M132").Value = -7121.804900000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 61374.082
$ws.Range("I62").Value = 4123.5
$ws.Range("J62").Value = 89999.375
$ws.Range("K62").Value = 4123.5
$ws.Range("L62").Value = 89999.375
$ws.Range("M62").Value = -3499.5
$ws.Range("N62").Value = -91247.375
$ws.Range("H65").Value = 61374.082
$ws.Range("I65").Value = 4123.5
$ws.Range("J65").Value = 89999.375
$ws.Range("K65").Value = 20617.5
$ws.Range("L65").Value = 449996.875
$ws.Range("M65").Value = -17497.5
$ws.Range("N65").Value = -456236.875
$ws.Range("H88").Value = 32535.5
$ws.Range("I88").Value = 50171
$ws.Range("J88").Value = 14900
$ws.Range("K88").Value = 50171
$ws.Range("L88").Value = 14900
$ws.Range("M88").Value = -49765
$ws.Range("N88").Value = -15712
$ws.Range("H91").Value = 32535.5
$ws.Range("I91").Value = 50171
$ws.Range("J91").Value = 14900
$ws.Range("K91").Value = 50171
$ws.Range("L91").Value = 14900
$ws.Range("M91").Value = -48767
$ws.Range("N91").Value = -17708
$ws.Range("H113").Value = 775.5
$ws.Range("I113").Value = 666.1429000000001
$ws.Range("K113").Value = 1998.4287
$ws.Range("M113").Value = 171.5712999999998
$ws.Range("H132").Value = 7618.737
$ws.Range("I132").Value = 9754.385
$ws.Range("J132").Value = 2991.5
$ws.Range("K132").Value = 29263.155
$ws.Range("L132").Value = 8974.5
$ws.Range("M132").Value = -26733.155
$ws.Range("N132").Value = -14034.5
